$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap columns B and C ---
# Was: B1=LNBSF00, C1=Date  -> Now: B1=Date, C1=LNBSF00
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "LNBSF00"

# --- Data rows 2-9: swap B (price) and C (date) values ---
# Column B becomes the Date column (numeric date serials, formatted),
# column C becomes the plain numeric (price) column (no special format).
$dates  = @(45734, 45733, 45730, 45729, 45728, 45735, 45736, 45737)
$prices = @(790.4, 795.08, 792.8440000000001, 803.816, 802.724, 808.9640000000001, 823.3680000000001, 823.9400000000001)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $dates[$i]
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 3).Value = $prices[$i]
    $ws.Cells.Item($r, 3).Style = "Normal"
}

# --- New row 10 ---
$ws.Cells.Item(10, 1).Value = 806.651
$ws.Cells.Item(10, 2).Value = 45740
$ws.Cells.Item(10, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(10, 3).Value = 788.6319999999999
$ws.Cells.Item(10, 3).Style = "Normal"
